# Apply scheduled-runner price/profit updates to the Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J17").Value = 5455937.5
$ws.Range("L17").Value = 16367812.5
$ws.Range("N17").Value = -16368148.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2624.2778
$ws.Range("I2").Value = 2472.7942
$ws.Range("K2").Value = 2472.7942
$ws.Range("M2").Value = -2359.7942

$ws.Range("H45").Value = 13126.818
$ws.Range("I45").Value = 16711.875
$ws.Range("K45").Value = 16711.875
$ws.Range("M45").Value = -16334.875

$ws.Range("H61").Value = 3991.5334
$ws.Range("I61").Value = 2997.2856
$ws.Range("K61").Value = 2997.2856
$ws.Range("M61").Value = -2785.2856

$ws.Range("H74").Value = 3666.2415
$ws.Range("I74").Value = 3165.4375
$ws.Range("J74").Value = 4282.615
$ws.Range("K74").Value = 3165.4375
$ws.Range("L74").Value = 4282.615
$ws.Range("M74").Value = -2291.4375
$ws.Range("N74").Value = -6030.615

$ws.Range("H77").Value = 3666.2415
$ws.Range("I77").Value = 3165.4375
$ws.Range("J77").Value = 4282.615
$ws.Range("K77").Value = 15827.1875
$ws.Range("L77").Value = 21413.075
$ws.Range("M77").Value = -11459.1875
$ws.Range("N77").Value = -30149.075

$ws.Range("H116").Value = 2624.2778
$ws.Range("I116").Value = 2472.7942
$ws.Range("K116").Value = 2472.7942
$ws.Range("M116").Value = -178.7941999999998

$ws.Range("H122").Value = 1734.1904
$ws.Range("I122").Value = 1642.6666
$ws.Range("K122").Value = 4927.9998
$ws.Range("M122").Value = -2477.9998

$ws.Range("H136").Value = 3991.5334
$ws.Range("I136").Value = 2997.2856
$ws.Range("K136").Value = 8991.856800000001
$ws.Range("M136").Value = -6441.856800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2624.2778
$ws.Range("I3").Value = 2472.7942
$ws.Range("K3").Value = 2472.7942
$ws.Range("M3").Value = -2358.7942

$ws.Range("H134").Value = 15865.333
$ws.Range("I134").Value = 4512.4287
$ws.Range("K134").Value = 13537.2861
$ws.Range("M134").Value = -11002.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2552.4753
$ws.Range("I31").Value = 1133.7826
$ws.Range("K31").Value = 1133.7826
$ws.Range("M31").Value = -838.7826

$ws.Range("H34").Value = 2552.4753
$ws.Range("I34").Value = 1133.7826
$ws.Range("K34").Value = 1133.7826
$ws.Range("M34").Value = -931.7826

$ws.Range("H58").Value = 3722.4707
$ws.Range("I58").Value = 3884.5
$ws.Range("J58").Value = 2966.3333
$ws.Range("K58").Value = 3884.5
$ws.Range("L58").Value = 2966.3333
$ws.Range("M58").Value = -3681.5
$ws.Range("N58").Value = -3372.3333

$ws.Range("H107").Value = 1974.75
$ws.Range("J107").Value = 1974.75
$ws.Range("L107").Value = 1974.75
$ws.Range("N107").Value = -5814.75

$ws.Range("H120").Value = 35000
$ws.Range("J120").Value = 35000
$ws.Range("L120").Value = 35000
$ws.Range("N120").Value = -42258

$ws.Range("H136").Value = 3722.4707
$ws.Range("I136").Value = 3884.5
$ws.Range("J136").Value = 2966.3333
$ws.Range("K136").Value = 11653.5
$ws.Range("L136").Value = 8898.999899999999
$ws.Range("M136").Value = -9103.5
$ws.Range("N136").Value = -13998.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H52").Value = 16652.4
$ws.Range("I52").Value = 12130
$ws.Range("J52").Value = 17783
$ws.Range("K52").Value = 12130
$ws.Range("L52").Value = 17783
$ws.Range("M52").Value = -11871
$ws.Range("N52").Value = -18301

$ws.Range("H97").Value = 870.91895
$ws.Range("I97").Value = 704.28
$ws.Range("J97").Value = 1218.0834
$ws.Range("K97").Value = 704.28
$ws.Range("L97").Value = 1218.0834
$ws.Range("M97").Value = -208.28
$ws.Range("N97").Value = -2210.0834

$ws.Range("H98").Value = 30992.5
$ws.Range("J98").Value = 30992.5
$ws.Range("L98").Value = 30992.5
$ws.Range("N98").Value = -36982.5

$ws.Range("H102").Value = 2082.5454
$ws.Range("I102").Value = 2087.0557
$ws.Range("J102").Value = 2062.25
$ws.Range("K102").Value = 2087.0557
$ws.Range("L102").Value = 2062.25
$ws.Range("M102").Value = -465.0556999999999
$ws.Range("N102").Value = -5306.25

$ws.Range("H111").Value = 28500
$ws.Range("J111").Value = 28500
$ws.Range("L111").Value = 28500
$ws.Range("N111").Value = -34634

$ws.Range("H122").Value = 2784.3572
$ws.Range("I122").Value = 1698.1111
$ws.Range("K122").Value = 5094.3333
$ws.Range("M122").Value = -2644.3333

$ws.Range("H126").Value = 2394.2917
$ws.Range("I126").Value = 2559.1428
$ws.Range("J126").Value = 1240.3334
$ws.Range("K126").Value = 7677.428400000001
$ws.Range("L126").Value = 3721.0002
$ws.Range("M126").Value = -5207.428400000001
$ws.Range("N126").Value = -8661.0002

$ws.Range("H132").Value = 8785.645500000001
$ws.Range("I132").Value = 10433.76
$ws.Range("J132").Value = 1918.5
$ws.Range("K132").Value = 31301.28
$ws.Range("L132").Value = 5755.5
$ws.Range("M132").Value = -28771.28
$ws.Range("N132").Value = -10815.5

$ws.Range("H136").Value = 11566
$ws.Range("J136").Value = 11566
$ws.Range("L136").Value = 34698
$ws.Range("N136").Value = -39798

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8478.380999999999
$ws.Range("I7").Value = 9115.0625
$ws.Range("K7").Value = 9115.0625
$ws.Range("M7").Value = -9003.0625

$ws.Range("H93").Value = 1850.6552
$ws.Range("I93").Value = 1828.8
$ws.Range("J93").Value = 1987.25
$ws.Range("K93").Value = 1828.8
$ws.Range("L93").Value = 1987.25
$ws.Range("M93").Value = -580.8
$ws.Range("N93").Value = -4483.25

$ws.Range("H100").Value = 5797.5
$ws.Range("I100").Value = 4646.25
$ws.Range("K100").Value = 4646.25
$ws.Range("M100").Value = -4105.25

$ws.Range("H126").Value = 8478.380999999999
$ws.Range("I126").Value = 9115.0625
$ws.Range("K126").Value = 27345.1875
$ws.Range("M126").Value = -24875.1875

$ws.Range("H132").Value = 2894.5293
$ws.Range("I132").Value = 2526.1538
$ws.Range("J132").Value = 4091.75
$ws.Range("K132").Value = 7578.4614
$ws.Range("L132").Value = 12275.25
$ws.Range("M132").Value = -5048.4614
$ws.Range("N132").Value = -17335.25

$ws.Range("H136").Value = 2366.5557
$ws.Range("I136").Value = 2065.88
$ws.Range("K136").Value = 6197.64
$ws.Range("M136").Value = -3647.64

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 15841.028
$ws.Range("I132").Value = 14947.542
$ws.Range("J132").Value = 17790.455
$ws.Range("K132").Value = 44842.626
$ws.Range("L132").Value = 53371.36500000001
$ws.Range("M132").Value = -42312.626
$ws.Range("N132").Value = -58431.36500000001

$ws.Range("H136").Value = 960.4828
$ws.Range("I136").Value = 984.7406999999999
$ws.Range("K136").Value = 2954.2221
$ws.Range("M136").Value = -404.2221
